$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 19.88361
$ws.Range("H2").Value = 59.65083
$ws.Range("I2").Value = 0.08287076207598171
$ws.Range("J2").Value = 0.08287076207598171
$ws.Range("M2").Value = 11.81073566666667
$ws.Range("N2").Value = 35.432207
$ws.Range("O2").Value = 0.3076347070004043
$ws.Range("P2").Value = 0.3076347070004043
$ws.Range("Q2").Value = 234.84006180909
$ws.Range("R2").Value = 2113.56055628181
$ws.Range("S2").Value = 0.02549392261014485
$ws.Range("T2").Value = 0.02549392261014485
$ws.Range("G3").Value = 19.88361
$ws.Range("H3").Value = 59.65083
$ws.Range("I3").Value = 0.08287076207598171
$ws.Range("J3").Value = 0.08287076207598171
$ws.Range("O3").Value = 0.3244871420261927
$ws.Range("P3").Value = 0.3244871420261927
$ws.Range("Q3").Value = 247.70475747909
$ws.Range("R3").Value = 2229.34281731181
$ws.Range("S3").Value = 0.0268904967435679
$ws.Range("T3").Value = 0.0268904967435679
$ws.Range("G4").Value = 19.88361
$ws.Range("H4").Value = 59.65083
$ws.Range("I4").Value = 0.08287076207598171
$ws.Range("J4").Value = 0.08287076207598171
$ws.Range("M4").Value = 3.197710666666667
$ws.Range("N4").Value = 9.593132000000001
$ws.Range("O4").Value = 0.08329089836363292
$ws.Range("P4").Value = 0.0832908983636329
$ws.Range("Q4").Value = 63.58203178884001
$ws.Range("R4").Value = 572.23828609956
$ws.Range("S4").Value = 0.006902380221387398
$ws.Range("T4").Value = 0.006902380221387397
$ws.Range("G5").Value = 19.88361
$ws.Range("H5").Value = 59.65083
$ws.Range("I5").Value = 0.08287076207598171
$ws.Range("J5").Value = 0.08287076207598171
$ws.Range("M5").Value = 8.081220666666667
$ws.Range("N5").Value = 24.243662
$ws.Range("O5").Value = 0.2104918797744333
$ws.Range("P5").Value = 0.2104918797744333
$ws.Range("Q5").Value = 160.68384005994
$ws.Range("R5").Value = 1446.15456053946
$ws.Range("S5").Value = 0.01744362248771321
$ws.Range("T5").Value = 0.01744362248771321
$ws.Range("G6").Value = 19.88361
$ws.Range("H6").Value = 59.65083
$ws.Range("I6").Value = 0.08287076207598171
$ws.Range("J6").Value = 0.08287076207598171
$ws.Range("M6").Value = 2.844675333333333
$ws.Range("N6").Value = 8.534026000000001
$ws.Range("O6").Value = 0.07409537283533685
$ws.Range("P6").Value = 0.07409537283533686
$ws.Range("Q6").Value = 56.56241490462001
$ws.Range("R6").Value = 509.06173414158
$ws.Range("S6").Value = 0.006140340013168358
$ws.Range("T6").Value = 0.006140340013168359
$ws.Range("I7").Value = 0.5704506040615172
$ws.Range("J7").Value = 0.5704506040615173
$ws.Range("M7").Value = 11.81073566666667
$ws.Range("N7").Value = 35.432207
$ws.Range("O7").Value = 0.3076347070004043
$ws.Range("P7").Value = 0.3076347070004043
$ws.Range("Q7").Value = 1616.54909115004
$ws.Range("R7").Value = 14548.94182035036
$ws.Range("S7").Value = 0.1754904044386685
$ws.Range("T7").Value = 0.1754904044386685
$ws.Range("I8").Value = 0.5704506040615172
$ws.Range("J8").Value = 0.5704506040615173
$ws.Range("O8").Value = 0.3244871420261927
$ws.Range("P8").Value = 0.3244871420261927
$ws.Range("S8").Value = 0.1851038861790369
$ws.Range("T8").Value = 0.1851038861790369
$ws.Range("I9").Value = 0.5704506040615172
$ws.Range("J9").Value = 0.5704506040615173
$ws.Range("M9").Value = 3.197710666666667
$ws.Range("N9").Value = 9.593132000000001
$ws.Range("O9").Value = 0.08329089836363292
$ws.Range("P9").Value = 0.0832908983636329
$ws.Range("Q9").Value = 437.6743682910401
$ws.Range("R9").Value = 3939.06931461936
$ws.Range("S9").Value = 0.04751334328436083
$ws.Range("T9").Value = 0.04751334328436083
$ws.Range("I10").Value = 0.5704506040615172
$ws.Range("J10").Value = 0.5704506040615173
$ws.Range("M10").Value = 8.081220666666667
$ws.Range("N10").Value = 24.243662
$ws.Range("O10").Value = 0.2104918797744333
$ws.Range("P10").Value = 0.2104918797744333
$ws.Range("Q10").Value = 1106.08604686264
$ws.Range("R10").Value = 9954.774421763759
$ws.Range("S10").Value = 0.1200752199673698
$ws.Range("T10").Value = 0.1200752199673698
$ws.Range("I11").Value = 0.5704506040615172
$ws.Range("J11").Value = 0.5704506040615173
$ws.Range("M11").Value = 2.844675333333333
$ws.Range("N11").Value = 8.534026000000001
$ws.Range("O11").Value = 0.07409537283533685
$ws.Range("P11").Value = 0.07409537283533686
$ws.Range("Q11").Value = 389.35401269672
$ws.Range("R11").Value = 3504.18611427048
$ws.Range("S11").Value = 0.04226775019208123
$ws.Range("T11").Value = 0.04226775019208125
$ws.Range("G12").Value = 28.56702866666667
$ws.Range("H12").Value = 85.701086
$ws.Range("I12").Value = 0.1190614499003492
$ws.Range("J12").Value = 0.1190614499003492
$ws.Range("M12").Value = 11.81073566666667
$ws.Range("N12").Value = 35.432207
$ws.Range("O12").Value = 0.3076347070004043
$ws.Range("P12").Value = 0.3076347070004043
$ws.Range("Q12").Value = 337.3976243640891
$ws.Range("R12").Value = 3036.578619276802
$ws.Range("S12").Value = 0.03662743425513724
$ws.Range("T12").Value = 0.03662743425513724
$ws.Range("G13").Value = 28.56702866666667
$ws.Range("H13").Value = 85.701086
$ws.Range("I13").Value = 0.1190614499003492
$ws.Range("J13").Value = 0.1190614499003492
$ws.Range("O13").Value = 0.3244871420261927
$ws.Range("P13").Value = 0.3244871420261927
$ws.Range("Q13").Value = 355.8804919114224
$ws.Range("R13").Value = 3202.924427202802
$ws.Range("S13").Value = 0.03863390960365903
$ws.Range("T13").Value = 0.03863390960365903
$ws.Range("G14").Value = 28.56702866666667
$ws.Range("H14").Value = 85.701086
$ws.Range("I14").Value = 0.1190614499003492
$ws.Range("J14").Value = 0.1190614499003492
$ws.Range("M14").Value = 3.197710666666667
$ws.Range("N14").Value = 9.593132000000001
$ws.Range("O14").Value = 0.08329089836363292
$ws.Range("P14").Value = 0.0832908983636329
$ws.Range("Q14").Value = 91.34909228237247
$ws.Range("R14").Value = 822.1418305413521
$ws.Range("S14").Value = 0.009916735122676757
$ws.Range("T14").Value = 0.009916735122676755
$ws.Range("G15").Value = 28.56702866666667
$ws.Range("H15").Value = 85.701086
$ws.Range("I15").Value = 0.1190614499003492
$ws.Range("J15").Value = 0.1190614499003492
$ws.Range("M15").Value = 8.081220666666667
$ws.Range("N15").Value = 24.243662
$ws.Range("O15").Value = 0.2104918797744333
$ws.Range("P15").Value = 0.2104918797744333
$ws.Range("Q15").Value = 230.8564624463258
$ws.Range("R15").Value = 2077.708162016932
$ws.Range("S15").Value = 0.02506146839819402
$ws.Range("T15").Value = 0.02506146839819402
$ws.Range("G16").Value = 28.56702866666667
$ws.Range("H16").Value = 85.701086
$ws.Range("I16").Value = 0.1190614499003492
$ws.Range("J16").Value = 0.1190614499003492
$ws.Range("M16").Value = 2.844675333333333
$ws.Range("N16").Value = 8.534026000000001
$ws.Range("O16").Value = 0.07409537283533685
$ws.Range("P16").Value = 0.07409537283533686
$ws.Range("Q16").Value = 81.26392179469291
$ws.Range("R16").Value = 731.375296152236
$ws.Range("S16").Value = 0.008821902520682153
$ws.Range("T16").Value = 0.008821902520682155
$ws.Range("G17").Value = 1.897401333333333
$ws.Range("H17").Value = 5.692203999999999
$ws.Range("I17").Value = 0.007907975184451771
$ws.Range("J17").Value = 0.007907975184451773
$ws.Range("M17").Value = 11.81073566666667
$ws.Range("N17").Value = 35.432207
$ws.Range("O17").Value = 0.3076347070004043
$ws.Range("P17").Value = 0.3076347070004043
$ws.Range("Q17").Value = 22.40970560158089
$ws.Range("R17").Value = 201.687350414228
$ws.Range("S17").Value = 0.002432767628835289
$ws.Range("T17").Value = 0.002432767628835289
$ws.Range("G18").Value = 1.897401333333333
$ws.Range("H18").Value = 5.692203999999999
$ws.Range("I18").Value = 0.007907975184451771
$ws.Range("J18").Value = 0.007907975184451773
$ws.Range("O18").Value = 0.3244871420261927
$ws.Range("P18").Value = 0.3244871420261927
$ws.Range("Q18").Value = 23.63732426424755
$ws.Range("R18").Value = 212.735918378228
$ws.Range("S18").Value = 0.002566036266816809
$ws.Range("T18").Value = 0.002566036266816809
$ws.Range("G19").Value = 1.897401333333333
$ws.Range("H19").Value = 5.692203999999999
$ws.Range("I19").Value = 0.007907975184451771
$ws.Range("J19").Value = 0.007907975184451773
$ws.Range("M19").Value = 3.197710666666667
$ws.Range("N19").Value = 9.593132000000001
$ws.Range("O19").Value = 0.08329089836363292
$ws.Range("P19").Value = 0.0832908983636329
$ws.Range("Q19").Value = 6.067340482547555
$ws.Range("R19").Value = 54.606064342928
$ws.Range("S19").Value = 0.0006586623573503037
$ws.Range("T19").Value = 0.0006586623573503037
$ws.Range("G20").Value = 1.897401333333333
$ws.Range("H20").Value = 5.692203999999999
$ws.Range("I20").Value = 0.007907975184451771
$ws.Range("J20").Value = 0.007907975184451773
$ws.Range("M20").Value = 8.081220666666667
$ws.Range("N20").Value = 24.243662
$ws.Range("O20").Value = 0.2104918797744333
$ws.Range("P20").Value = 0.2104918797744333
$ws.Range("Q20").Value = 15.33331886789422
$ws.Range("R20").Value = 137.999869811048
$ws.Range("S20").Value = 0.001664564561784824
$ws.Range("T20").Value = 0.001664564561784825
$ws.Range("G21").Value = 1.897401333333333
$ws.Range("H21").Value = 5.692203999999999
$ws.Range("I21").Value = 0.007907975184451771
$ws.Range("J21").Value = 0.007907975184451773
$ws.Range("M21").Value = 2.844675333333333
$ws.Range("N21").Value = 8.534026000000001
$ws.Range("O21").Value = 0.07409537283533685
$ws.Range("P21").Value = 0.07409537283533686
$ws.Range("Q21").Value = 5.397490770367111
$ws.Range("R21").Value = 48.577416933304
$ws.Range("S21").Value = 0.0005859443696645456
$ws.Range("T21").Value = 0.0005859443696645459
$ws.Range("G22").Value = 52.71596533333334
$ws.Range("H22").Value = 158.147896
$ws.Range("I22").Value = 0.2197092087777001
$ws.Range("J22").Value = 0.2197092087777001
$ws.Range("M22").Value = 11.81073566666667
$ws.Range("N22").Value = 35.432207
$ws.Range("O22").Value = 0.3076347070004043
$ws.Range("P22").Value = 0.3076347070004043
$ws.Range("Q22").Value = 622.6143319651636
$ws.Range("R22").Value = 5603.528987686472
$ws.Range("S22").Value = 0.06759017806761844
$ws.Range("T22").Value = 0.06759017806761844
$ws.Range("G23").Value = 52.71596533333334
$ws.Range("H23").Value = 158.147896
$ws.Range("I23").Value = 0.2197092087777001
$ws.Range("J23").Value = 0.2197092087777001
$ws.Range("O23").Value = 0.3244871420261927
$ws.Range("P23").Value = 0.3244871420261927
$ws.Range("Q23").Value = 656.7215615358302
$ws.Range("R23").Value = 5910.494053822472
$ws.Range("S23").Value = 0.071292813233112
$ws.Range("T23").Value = 0.071292813233112
$ws.Range("G24").Value = 52.71596533333334
$ws.Range("H24").Value = 158.147896
$ws.Range("I24").Value = 0.2197092087777001
$ws.Range("J24").Value = 0.2197092087777001
$ws.Range("M24").Value = 3.197710666666667
$ws.Range("N24").Value = 9.593132000000001
$ws.Range("O24").Value = 0.08329089836363292
$ws.Range("P24").Value = 0.0832908983636329
$ws.Range("Q24").Value = 168.5704046500302
$ws.Range("R24").Value = 1517.133641850272
$ws.Range("S24").Value = 0.01829977737785763
$ws.Range("T24").Value = 0.01829977737785762
$ws.Range("G25").Value = 52.71596533333334
$ws.Range("H25").Value = 158.147896
$ws.Range("I25").Value = 0.2197092087777001
$ws.Range("J25").Value = 0.2197092087777001
$ws.Range("M25").Value = 8.081220666666667
$ws.Range("N25").Value = 24.243662
$ws.Range("O25").Value = 0.2104918797744333
$ws.Range("P25").Value = 0.2104918797744333
$ws.Range("Q25").Value = 426.0093485150169
$ws.Range("R25").Value = 3834.084136635152
$ws.Range("S25").Value = 0.04624700435937153
$ws.Range("T25").Value = 0.04624700435937153
$ws.Range("G26").Value = 52.71596533333334
$ws.Range("H26").Value = 158.147896
$ws.Range("I26").Value = 0.2197092087777001
$ws.Range("J26").Value = 0.2197092087777001
$ws.Range("M26").Value = 2.844675333333333
$ws.Range("N26").Value = 8.534026000000001
$ws.Range("O26").Value = 0.07409537283533685
$ws.Range("P26").Value = 0.07409537283533686
$ws.Range("Q26").Value = 149.9598062565885
$ws.Range("R26").Value = 1349.638256309296
$ws.Range("S26").Value = 0.01627943573974056
$ws.Range("T26").Value = 0.01627943573974056
